$d = $word.ActiveDocument

# The paragraph originally read:
#   "9号开始第一天去部门报到，"
# It is edited so that "报到，" is dropped and replaced by a much longer
# continuation, giving the final text:
#   "9号开始第一天去部门就差点迟到，因为不认识路，找工位费了点时间。
#    当时还是在E2-5F-1岛，第一次晨会的时候听到大家讨论的内容感觉既兴奋又
#    紧张。兴奋是因为终于加入了紧贴科技最前沿的手机行业，作为数码迷的我
#    喜不自胜。紧张是因为本次跳槽跨行业，对驱动知之甚少，担心跟不上节奏。"
# The inserted text keeps exactly the same run formatting (苹方-简, bold,
# sz 28 / szCs 36) as the run it continues, so a single Find/Execute
# replacement (which preserves the formatting of the matched run) gives
# the correct resulting document content.

$oldText = "9号开始第一天去部门报到，"
$newText = "9号开始第一天去部门就差点迟到，因为不认识路，找工位费了点时间。当时还是在E2-5F-1岛，第一次晨会的时候听到大家讨论的内容感觉既兴奋又紧张。兴奋是因为终于加入了紧贴科技最前沿的手机行业，作为数码迷的我喜不自胜。紧张是因为本次跳槽跨行业，对驱动知之甚少，担心跟不上节奏。"

$rng = $d.Content
$rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
